$d = $word.ActiveDocument

# --- Repository URL: the 3 existing runs already concatenate to the correct
# final text (mlops_st + a + rter = mlops_starter); merge them into a single
# run via a self-replace, then restore the Hyperlink character style that
# Find/Replace drops on the merged run.
$found = $d.Content.Find.Execute("https://github.com/Mitch1789/mlops_starter", $true, $false, $false, $false, $false, $true, 1, $false, "https://github.com/Mitch1789/mlops_starter", 2)
if (-not $found) { Write-Output "MISS: url-merge" }
$r = $d.Content
$found = $r.Find.Execute("https://github.com/Mitch1789/mlops_starter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "MISS: url-restyle" }
$r.Style = "Hyperlink"

# --- Text edits. Each Find/Replace is scoped to the minimal changed substring
# (rather than whole-sentence) so untouched punctuation/quotes elsewhere in the
# sentence are never re-typed by the replace (which would trigger Word's
# smart-quote autocorrect and corrupt straight apostrophes into curly ones).
$found = $d.Content.Find.Execute("Starter — Bank", $true, $false, $false, $false, $false, $true, 1, $false, "Starter Bank", 2)
if (-not $found) { Write-Output "MISS: title-emdash" }
$found = $d.Content.Find.Execute("August 10,", $true, $false, $false, $false, $false, $true, 1, $false, "August 13,", 2)
if (-not $found) { Write-Output "MISS: date" }
$found = $d.Content.Find.Execute("Dataset — ", $true, $false, $false, $false, $false, $true, 1, $false, "Dataset: ", 2)
if (-not $found) { Write-Output "MISS: dataset" }
$found = $d.Content.Find.Execute("Versioning — ", $true, $false, $false, $false, $false, $true, 1, $false, "Versioning: ", 2)
if (-not $found) { Write-Output "MISS: versioning" }
$found = $d.Content.Find.Execute("Pipeline — Four", $true, $false, $false, $false, $false, $true, 1, $false, "Pipeline: Four", 2)
if (-not $found) { Write-Output "MISS: pipeline" }
$found = $d.Content.Find.Execute("Baseline Model — ", $true, $false, $false, $false, $false, $true, 1, $false, "Baseline Model: ", 2)
if (-not $found) { Write-Output "MISS: baseline-model" }
$found = $d.Content.Find.Execute("Inference — ", $true, $false, $false, $false, $false, $true, 1, $false, "Inference: ", 2)
if (-not $found) { Write-Output "MISS: inference" }
$found = $d.Content.Find.Execute("Containerization — ", $true, $false, $false, $false, $false, $true, 1, $false, "Containerization: ", 2)
if (-not $found) { Write-Output "MISS: containerization" }
$found = $d.Content.Find.Execute("CI/CD — ", $true, $false, $false, $false, $false, $true, 1, $false, "CI/CD - ", 2)
if (-not $found) { Write-Output "MISS: ci-cd" }
$found = $d.Content.Find.Execute("Monitoring — ", $true, $false, $false, $false, $false, $true, 1, $false, "Monitoring - ", 2)
if (-not $found) { Write-Output "MISS: monitoring" }
$found = $d.Content.Find.Execute("Class 1 — ", $true, $false, $false, $false, $false, $true, 1, $false, "Class 1 - ", 2)
if (-not $found) { Write-Output "MISS: class1" }
$found = $d.Content.Find.Execute("System metrics — ", $true, $false, $false, $false, $false, $true, 1, $false, "System metrics - ", 2)
if (-not $found) { Write-Output "MISS: system-metrics" }
$found = $d.Content.Find.Execute("Data/Concept drift — ", $true, $false, $false, $false, $false, $true, 1, $false, "Data/Concept drift - ", 2)
if (-not $found) { Write-Output "MISS: data-concept-drift" }
$found = $d.Content.Find.Execute("Logging — ", $true, $false, $false, $false, $false, $true, 1, $false, "Logging - ", 2)
if (-not $found) { Write-Output "MISS: logging" }
$found = $d.Content.Find.Execute("Container parity matters — ", $true, $false, $false, $false, $false, $true, 1, $false, "Container parity matters - ", 2)
if (-not $found) { Write-Output "MISS: container-parity" }
$found = $d.Content.Find.Execute("Windows hygiene — ", $true, $false, $false, $false, $false, $true, 1, $false, "Windows hygiene - ", 2)
if (-not $found) { Write-Output "MISS: windows-hygiene" }
$found = $d.Content.Find.Execute("Cloud consistency — ", $true, $false, $false, $false, $false, $true, 1, $false, "Cloud consistency - ", 2)
if (-not $found) { Write-Output "MISS: cloud-consistency" }

Write-Output "DONE"
